# The commit reshuffles the observation rows 3-18 of the sheet: each
# destination row ends up holding the full set of field values that used
# to live on a different source row (the upstream export got re-sorted /
# re-paginated), while the column layout (A:AY) itself is untouched.
#
# Strategy:
#   1. Snapshot every involved source row's full value vector (A:AY)
#      before making any writes, so overlapping reads/writes can't
#      clobber each other while we shuffle things around.
#   2. Columns Y and AA hold plain "YYYY-MM-DD" text (not real Excel
#      dates) in the original file. Writing such a string back through
#      .Value would normally get auto-coerced into a date serial, same
#      as typing it into Excel by hand, so those two columns are
#      temporarily forced to Text format for the duration of the write
#      and then restored to the default style afterwards.
#   3. Write each destination row using the snapshot captured in step 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = "A"
$lastCol  = "AY"
$firstRow = 3
$lastRow  = 18

# new row number -> old row number that its data comes from
$mapping = @{
    3  = 5
    4  = 8
    5  = 18
    6  = 3
    7  = 4
    8  = 6
    9  = 7
    10 = 9
    11 = 10
    12 = 11
    13 = 12
    14 = 13
    15 = 14
    16 = 15
    17 = 16
    18 = 17
}

$dateColRange1 = "Y$firstRow" + ":" + "Y$lastRow"
$dateColRange2 = "AA$firstRow" + ":" + "AA$lastRow"

# Avoid Excel's automatic text->date coercion for the plain-text date
# strings while we copy them around.
$ws.Range($dateColRange1).NumberFormat = "@"
$ws.Range($dateColRange2).NumberFormat = "@"

# 1) Snapshot every involved source row before any writes happen.
$snapshot = @{}
foreach ($r in $firstRow..$lastRow) {
    $rangeAddr = "$firstCol$r" + ":" + "$lastCol$r"
    $snapshot[$r] = $ws.Range($rangeAddr).Value()
}

# 2) Write each destination row using the snapshot captured above.
foreach ($newRow in $firstRow..$lastRow) {
    $oldRow = $mapping[$newRow]
    $rangeAddr = "$firstCol$newRow" + ":" + "$lastCol$newRow"
    $ws.Range($rangeAddr).Value = $snapshot[$oldRow]
}

# Restore the default (un-formatted) style now that the text values are
# safely locked in, matching the original file's formatting.
$ws.Range($dateColRange1).Style = "Normal"
$ws.Range($dateColRange2).Style = "Normal"
